# Natmi following Dr Hou advice
# Adds the "ECs" sending/target cluster to the Efna5 -> Ephb6 LR-pair sheet,
# recomputing the full 3x3 cluster cross-join (ECs, FAPs, sCs) that replaces
# the previous 2x2 (FAPs, sCs) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newData = @(
    @("ECs", "Efna5", "Ephb6", "ECs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 2, 0.6666666666666666, 0.2721246666666666, 0.8163739999999999, 0.1154759865526449, 0.1154759865526449, 0.06282133992333333, 0.5653920593099999, 0.007364255951776857, 0.007364255951776857),
    @("ECs", "Efna5", "Ephb6", "FAPs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 3, 1, 0.8319233333333332, 2.49577, 0.3530263126440755, 0.3530263126440755, 0.1920536611166666, 1.72848295005, 0.02251356495523636, 0.02251356495523636),
    @("ECs", "Efna5", "Ephb6", "sCs", 2, 0.6666666666666666, 0.230855, 0.692565, 0.06377305075821572, 0.06377305075821572, 3, 1, 1.252499666666667, 3.757499, 0.5314977008032796, 0.5314977008032796, 0.2891458105483334, 2.602312294935, 0.0338952298512025, 0.0338952298512025),
    @("FAPs", "Efna5", "Ephb6", "ECs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 2, 0.6666666666666666, 0.2721246666666666, 0.8163739999999999, 0.1154759865526449, 0.1154759865526449, 0.7893605471728887, 7.104244924556, 0.09253309647183565, 0.09253309647183566),
    @("FAPs", "Efna5", "Ephb6", "FAPs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 3, 1, 0.8319233333333332, 2.49577, 0.3530263126440755, 0.3530263126440755, 2.413186079931111, 21.71867471938, 0.2828866747122192, 0.2828866747122193),
    @("FAPs", "Efna5", "Ephb6", "sCs", 3, 1, 2.900731333333333, 8.702194, 0.8013189515350044, 0.8013189515350045, 3, 1, 1.252499666666667, 3.757499, 0.5314977008032796, 0.5314977008032796, 3.633165028089556, 32.698485252806, 0.4258991803509495, 0.4258991803509496),
    @("sCs", "Efna5", "Ephb6", "ECs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 2, 0.6666666666666666, 0.2721246666666666, 0.8163739999999999, 0.1154759865526449, 0.1154759865526449, 0.1328947115051111, 1.196052403546, 0.01557863412903235, 0.01557863412903235),
    @("sCs", "Efna5", "Ephb6", "FAPs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 3, 1, 0.8319233333333332, 2.49577, 0.3530263126440755, 0.3530263126440755, 0.4062778017588889, 3.65650021583, 0.04762607297661986, 0.04762607297661986),
    @("sCs", "Efna5", "Ephb6", "sCs", 3, 1, 0.4883596666666667, 1.465079, 0.1349079977067798, 0.1349079977067798, 3, 1, 1.252499666666667, 3.757499, 0.5314977008032796, 0.5314977008032796, 0.6116703197134445, 5.505032877421001, 0.07170329060112758, 0.07170329060112758)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $row = $newData[$i]
    $r = $startRow + $i
    for ($c = 1; $c -le $row.Count; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
